# Re-style the three tables that used the deck's custom default table
# style ("Table_0", GUID A263DC0B-...) so that they use PowerPoint's
# built-in "Medium Style 2 - Accent 1" table style
# (GUID 0FBBE7AE-6FFD-4052-8BA6-E1CBA36C2E7A) instead — matching a
# Table Styles gallery selection made on each of the three tables.

$p = $ppt.ActivePresentation

$oldStyleId = "{A263DC0B-0AC0-4D2C-B9D6-BB7C4A4E65FF}"
$newStyleId = "{0FBBE7AE-6FFD-4052-8BA6-E1CBA36C2E7A}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.HasTable) {
            $table = $shape.Table

            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
